$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 2
$ws.Range("H2").Value = 3
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 5
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 5.5
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 6.5
$ws.Range("AE2").Value = 19
$ws.Range("AJ2").Value = 17
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 81
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 151

# Remove the second match (row 3) entirely
$ws.Rows.Item(3).Delete()
